$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cosmetic re-save touches (Hungarian Excel -> English Excel locale):
# builtin "Normal" cell style and default theme display name.
$wb.Styles.Item(1).Name = "Normal"
$wb.Theme.Name = "Office Theme"

# New bestiary entries (name only - stats to be filled in later).
# NOTE: "Specter" (row 34) is written before "Succubus" (row 33) on
# purpose, so the shared-string table picks up the same index order
# as the source workbook (Specter before Succubus) even though the
# visible row order places Succubus above Specter.
$ws.Range("A18").Value = "Skeleton"
$ws.Range("A19").Value = "Ent"
$ws.Range("A20").Value = "Kobold"
$ws.Range("A21").Value = "Hag"
$ws.Range("A22").Value = "Wyvern"
$ws.Range("A23").Value = "Griffon"
$ws.Range("A24").Value = "Manticore"
$ws.Range("A25").Value = "Centaur"
$ws.Range("A26").Value = "Carnivorous Plant"
$ws.Range("A27").Value = "Mad Knight"
$ws.Range("A28").Value = "Vampire"
$ws.Range("A29").Value = "Gargoyle"
$ws.Range("A30").Value = "Cannibal"
$ws.Range("A31").Value = "Basilisk"
$ws.Range("A32").Value = "Acromantula"
$ws.Range("A34").Value = "Specter"
$ws.Range("A33").Value = "Succubus"
$ws.Range("A35").Value = "Toxic Spore"
$ws.Range("A36").Value = "Fungoid"
$ws.Range("A37").Value = "Shadow Demon"
$ws.Range("A38").Value = "Infernal Fiend"
$ws.Range("A39").Value = "Naga"
$ws.Range("A40").Value = "Sea Serpent"
$ws.Range("A41").Value = "The Krakken"

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor / selection on the last-added row, matching the
# author's final cell selection.
[void]$ws.Range("A33").Select()
